# Fixed a bug in stats2
# The underlying test data rows (A2:F21) were regenerated/reordered; this
# script rewrites those rows to the corrected values while leaving the
# header row, the trailing per-reel summary rows (22:25) and the totals
# row (26) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(401,9,48,67,75,45),
    @(701,3,90,45,97,15),
    @(201,9,30,15,45,30),
    @(801,3,67,65,52,45),
    @(101,9,30,15,60,15),
    @(1203,3,15,15,15,15),
    @(501,9,52,30,75,45),
    @(601,9,60,67,60,42),
    @(1201,2,10,10,10,10),
    @(1202,2,10,10,10,10),
    @(901,16,15,45,60,60),
    @(902,1,0,0,0,0),
    @(1001,18,30,75,60,72),
    @(301,6,45,30,60,45),
    @(1,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(502,0,4,0,0,0),
    @(2,0,2,2,2,2),
    @(1101,0,15,30,30,0),
    @(802,0,4,5,4,0)
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $values[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
